$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$newName = "1013-MS-EI-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-LateRepayment"

# Update the product name value on both sheets (space removed from "Late Repayment").
$ws1.Range("B1").Value = $newName
$ws2.Range("B1").Value = $newName

# Update selection on ProductLoanOutput (not the active tab any more) to B1.
$ws2.Activate()
$ws2.Range("B1").Select()

# ProductLoanInput becomes the active/selected tab, scrolled to top with B1 selected.
$ws1.Activate()
$ws1.Range("B1").Select()
